$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cells, matching the formatting (bold, border, centered)
# of the existing header row by copying the format from the last header cell.
$xlPasteFormats = -4122
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}
